$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.5
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("W2").Value = 19
$ws.Range("X2").Value = 29
$ws.Range("AA2").Value = 34
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 8
$ws.Range("AH2").Value = 10
$ws.Range("AR2").Value = 81
$ws.Range("AT2").Value = 3.5
$ws.Range("K3").Value = 2.75
$ws.Range("N3").Value = 15
$ws.Range("AH3").Value = 29
$ws.Range("AJ3").Value = 34
$ws.Range("AK3").Value = 151
$ws.Range("AL3").Value = 81
$ws.Range("AM3").Value = 67
$ws.Range("AQ3").Value = 13
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 4
$ws.Range("S4").Value = 1.33
$ws.Range("T4").Value = 3.25
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 2.1
$ws.Range("Z4").Value = 19
$ws.Range("AA4").Value = 15
$ws.Range("AI4").Value = 19
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 26
$ws.Range("AO4").Value = 11
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 34
$ws.Range("AT4").Value = 3.25
$ws.Range("AY4").Value = 19
$ws.Range("BD5").Value = 151
$ws.Range("I6").Value = 3.6
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.98
$ws.Range("AA6").Value = 17
$ws.Range("AG6").Value = 201
$ws.Range("AS6").Value = 151
$ws.Range("AW6").Value = 126
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 2.3
$ws.Range("J8").Value = 3.5
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("Q8").Value = 1.73
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.2
$ws.Range("W8").Value = 11
$ws.Range("AB8").Value = 29
$ws.Range("AE8").Value = 13
$ws.Range("AG8").Value = 151
$ws.Range("AH8").Value = 9.5
$ws.Range("AI8").Value = 12
$ws.Range("AJ8").Value = 9
$ws.Range("AK8").Value = 21
$ws.Range("AP8").Value = 23
$ws.Range("AS8").Value = 151
$ws.Range("AT8").Value = 3
$ws.Range("AV8").Value = 51
$ws.Range("AW8").Value = 501
$ws.Range("I11").Value = 4
$ws.Range("L11").Value = 4.33
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("Q11").Value = 1.8
$ws.Range("R11").Value = 2
$ws.Range("U11").Value = 1.67
$ws.Range("V11").Value = 2.1
$ws.Range("X11").Value = 9.5
$ws.Range("AC11").Value = 13
$ws.Range("AE11").Value = 13
$ws.Range("AL11").Value = 29
$ws.Range("AQ11").Value = 29
$ws.Range("AU11").Value = 7.5
$ws.Range("AZ11").Value = 26
$ws.Range("BC11").Value = 151
$ws.Range("G12").Value = 1.65
$ws.Range("H12").Value = 3.55
$ws.Range("I12").Value = 4.9
$ws.Range("J12").Value = 2.2
$ws.Range("K12").Value = 2.12
$ws.Range("O12").Value = 1.27
$ws.Range("P12").Value = 3.1
$ws.Range("Q12").Value = 1.78
$ws.Range("R12").Value = 1.82
$ws.Range("W12").Value = 7
$ws.Range("X12").Value = 7.9
$ws.Range("Y12").Value = 7.9
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 13
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 7
$ws.Range("AF12").Value = 75
$ws.Range("AG12").Value = 600
$ws.Range("AH12").Value = 13
$ws.Range("AI12").Value = 29
$ws.Range("AO12").Value = 8
$ws.Range("AQ12").Value = 26
$ws.Range("AT12").Value = 2.55
$ws.Range("AV12").Value = 70
$ws.Range("AX12").Value = 6.5
$ws.Range("AY12").Value = 29
$ws.Range("AZ12").Value = 35
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 8.9
$ws.Range("G14").Value = 2.35
$ws.Range("H14").Value = 3.15
$ws.Range("I14").Value = 2.8
$ws.Range("L14").Value = 3.35
$ws.Range("N14").Value = 7.9
$ws.Range("V14").Value = 2.32
$ws.Range("W14").Value = 10
$ws.Range("AC14").Value = 7.9
$ws.Range("AD14").Value = 6.4
$ws.Range("AE14").Value = 11.25
$ws.Range("AH14").Value = 10.75
$ws.Range("AI14").Value = 16.5
$ws.Range("AM14").Value = 25
$ws.Range("AU14").Value = 6.3
$ws.Range("AX14").Value = 5
$ws.Range("AY14").Value = 15.5
$ws.Range("AZ14").Value = 19.5
